$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price/volume values for existing rows (rows 2-46 keep same coins,
# only the Price (D) and Volume(1h) (E) figures are refreshed).
# Rows 47-51 shift down by one: a new coin (RenzoRestakedETH) is inserted
# before the old "Stellar" row, the following rows (Stellar, Cosmos, ONDO,
# InjectiveProtocol) move down one slot with updated figures, and the former
# last row (SuiNetwork) drops off the bottom of the table.
#
# Price values that look like plain numbers (e.g. "537.42") must stay plain
# text cells (matching the source data, which stores every Price/Volume
# figure as text) instead of being auto-converted to numbers by Excel, so
# those assignments are wrapped with a Text number format that is cleared
# again right after the value is set (keeps the cell's style index at its
# original default while forcing a text value).

$ws.Range("D2").Value = "58.532.71"
$ws.Range("E2").Value = "  +2.73%  "
$ws.Range("D3").Value = "3.157.94"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.88"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +10.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.33"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.422"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.139"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("D13").Value = "3.708.99"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.15"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000169"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.68%  "
$ws.Range("D16").Value = "58.594.40"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("D17").Value = "3.156.71"
$ws.Range("E17").Value = "  +3.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.22"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.06"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.25"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "378.19"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +9.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.76"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.993"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.42"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.516"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.167"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.981"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.10"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +14.13%  "
$ws.Range("D29").Value = "0.0₃0873"
$ws.Range("E29").Value = "  +4.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.18"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +7.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.89"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.90"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.17"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +8.01%  "
$ws.Range("E34").Value = "  +6.37%  "
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.25"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.37"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +13.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.48"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.67"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.51%  "
$ws.Range("D40").Value = "2.634.67"
$ws.Range("E40").Value = "  +9.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0682"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.22"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.67"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.706"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("B47").Value = "RenzoRestakedETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D47").Value = "3.206.62"
$ws.Range("E47").Value = "  +3.54%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.103"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +12.83%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.22"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.65%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.981"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.49%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.44%  "
